$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.706.78'
$ws.Range("E2").Value = '  -0.21%  '

$ws.Range("D3").Value = '3.525.39'
$ws.Range("E3").Value = '  -1.04%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '614.70'
$ws.Range("E5").Value = '  +0.07%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.48'
$ws.Range("E6").Value = '  +0.59%  '

$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.610'
$ws.Range("E7").Value = '  -1.44%  '

$ws.Range("B8").Value = 'LidoStakedEther'
$ws.Range("C8").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D8").Value = '3.517.64'
$ws.Range("E8").Value = '  -1.15%  '

$ws.Range("E9").Value = '  -0.03%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.197'
$ws.Range("E10").Value = '  -0.56%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.28'
$ws.Range("E11").Value = '  +3.04%  '

$ws.Range("E12").Value = '  +0.38%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '46.53'
$ws.Range("E13").Value = '  -0.71%  '

$ws.Range("E14").Value = '  -0.85%  '

$ws.Range("D15").Value = '4.095.05'
$ws.Range("E15").Value = '  -1.01%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '8.43'
$ws.Range("E16").Value = '  +0.22%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '617.57'
$ws.Range("E17").Value = '  -0.42%  '

$ws.Range("D18").Value = '3.524.26'
$ws.Range("E18").Value = '  -1.00%  '

$ws.Range("D19").Value = '70.741.38'
$ws.Range("E19").Value = '  -0.31%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.121'
$ws.Range("E20").Value = '  +1.46%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.75'
$ws.Range("E21").Value = '  +1.97%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.885'
$ws.Range("E22").Value = '  +0.08%  '

$ws.Range("E23").Value = '  -5.63%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '15.72'
$ws.Range("E24").Value = '  -0.31%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '98.48'
$ws.Range("E25").Value = '  +1.68%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.79'
$ws.Range("E26").Value = '  -0.98%  '

$ws.Range("E27").Value = '  -0.06%  '

$ws.Range("E28").Value = '  -0.38%  '

$ws.Range("E29").Value = '  +0.96%  '

$ws.Range("E30").Value = '  +0.59%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.05'
$ws.Range("E31").Value = '  -1.07%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.16'
$ws.Range("E32").Value = '  -4.35%  '

$ws.Range("E33").Value = '  -0.16%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.85'
$ws.Range("E34").Value = '  -2.47%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '620.12'
$ws.Range("E35").Value = '  +8.01%  '

$ws.Range("E36").Value = '  -0.96%  '

$ws.Range("E37").Value = '  +0.04%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0484'
$ws.Range("E38").Value = '  +2.49%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.50'
$ws.Range("E39").Value = '  -3.83%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '57.07'
$ws.Range("E40").Value = '  -1.31%  '

$ws.Range("E41").Value = '  -0.03%  '

$ws.Range("E42").Value = '  +1.45%  '

$ws.Range("D43").Value = '3.373.06'
$ws.Range("E43").Value = '  +0.28%  '

$ws.Range("D44").Value = '0.0₃0739'
$ws.Range("E44").Value = '  +3.87%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.313'
$ws.Range("E45").Value = '  -2.86%  '

$ws.Range("E46").Value = '  -3.64%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '32.28'
$ws.Range("E47").Value = '  -2.63%  '

$ws.Range("E48").Value = '  -2.63%  '

$ws.Range("E49").Value = '  +0.32%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.67'
$ws.Range("E50").Value = '  -0.34%  '

$ws.Range("E51").Value = '  -0.01%  '

